$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.451.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "'1.828.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "'315.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("E7").Value = "  -4.00%  "
$ws.Range("D8").Value = "'0.3924"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "'0.07666"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.61%  "
$ws.Range("D10").Value = "'41.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "'1.110"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "'21.07"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").Value = "'6.307"
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "'7.530"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.04%  "
$ws.Range("D16").Value = "'1.823.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E17").Value = "  +4.63%  "
$ws.Range("D18").Value = "'0.00001100"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.03%  "
$ws.Range("D19").Value = "'0.06696"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.73%  "
$ws.Range("D20").Value = "'17.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.59%  "
$ws.Range("E21").Value = "  +0.42%  "
$ws.Range("D22").Value = "'6.143"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").Value = "'28.491.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'2.256"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.78%  "
$ws.Range("D26").Value = "'20.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").Value = "'156.97"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "'2.039.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("D29").Value = "'2.401"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("D30").Value = "'124.55"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "'1.114"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("D32").Value = "'0.1083"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").Value = "'5.664"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("D34").Value = "'3.666"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").Value = "'0.07013"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.99%  "
$ws.Range("D36").Value = "'0.2213"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("D37").Value = "'8.940"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.30%  "
$ws.Range("D38").Value = "'0.02321"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").Value = "'5.161"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.96%  "
$ws.Range("D40").Value = "'0.6269"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("D41").Value = "'11.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.91%  "
$ws.Range("D42").Value = "'1.175"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").Value = "'1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("D45").Value = "'13.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("D46").Value = "'0.5893"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.88%  "
$ws.Range("D47").Value = "'3.713"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").Value = "'124.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").Value = "'1.196"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").Value = "'0.06928"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.76%  "
